# Commit: "In bilateral plot, color industries according to industry groups"
#
# The "data" sheet lists NAICS industry codes/names alongside a group id
# (column D) used elsewhere to color industries by their parent group.
# The sheet used to also contain four "group header" rows (placeholder
# codes R1-R4 in column A, paired with the group's own name in column B,
# e.g. "Chemicals", "Computer and Electronic Products",
# "Transportation Equipment", "Miscellaneous Manufacturing") that aren't
# needed now that industries are colored by group, so those rows are
# removed, shrinking the table from 30 to 26 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so earlier row numbers stay valid as later rows shift up.
$ws.Rows(28).Delete()   # R4 | Miscellaneous Manufacturing | 339
$ws.Rows(23).Delete()   # R3 | Transportation Equipment    | 336
$ws.Rows(16).Delete()   # R2 | Computer and Electronic Products | 334
$ws.Rows(6).Delete()    # R1 | Chemicals                    | 325
